# Swap the contents of columns A and B (both header row and all data rows).
# Before: A = 分类名称 (category, constant "水生根茎类"), B = 单品名称 (item name)
# After:  A = 单品名称 (item name),                       B = 分类名称 (category, constant "水生根茎类")

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 1; $r -le $lastRow; $r++) {
    $aVal = $ws.Cells.Item($r, 1).Value2
    $bVal = $ws.Cells.Item($r, 2).Value2
    $ws.Cells.Item($r, 1).Value2 = $bVal
    $ws.Cells.Item($r, 2).Value2 = $aVal
}
